$d = $word.ActiveDocument

# The big progress-tracking table is the 3rd table in the document.
$tbl = $d.Tables.Item(3)

# --- Row for Task 23 (date 20/11/17): Status cell ---
# "In Progress (database prototype to be done)" -> "Complete"
$statusCell23 = $tbl.Cell(24, 3)
$statusRange23 = $d.Range($statusCell23.Range.Start, $statusCell23.Range.End)
$statusRange23.Find.Execute("In Progress (database prototype to be done)", $true, $false, $false, $false, $false, $true, 0, $false, "Complete", 2)

# --- Row for Task 23 (date 20/11/17): Resources cell ---
# "E - 6h  C - 3h" -> "E - 8h  C - 3h" (just the "6" becomes "8")
$resCell23 = $tbl.Cell(24, 4)
$resRange23 = $d.Range($resCell23.Range.Start, $resCell23.Range.End)
$resRange23.Find.Execute("6h", $true, $false, $false, $false, $false, $true, 0, $false, "8h", 2)

# --- Row for Task 25 (date 15/11/17): Status cell ---
# "In progress" -> "Complete"
$statusCell25 = $tbl.Cell(26, 3)
$statusRange25 = $d.Range($statusCell25.Range.Start, $statusCell25.Range.End)
$statusRange25.Find.Execute("In progress", $true, $false, $false, $false, $false, $true, 0, $false, "Complete", 2)

# --- Row for Task 27 (date 21/11/17): Resources cell ---
# "E -2h" stays the same text, but the cursor was left between "E -2" and "h"
# (Word drops a "_GoBack" bookmark at the last edit position). Split the run
# there and (re)plant the _GoBack bookmark, which also removes it from its
# old location at the end of the document.
$resCell27 = $tbl.Cell(28, 4)
$resRange27 = $d.Range($resCell27.Range.Start, $resCell27.Range.End)
$resRange27.Find.Execute("E -2h", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$splitPoint = $resRange27.Start + 4
$goBack = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $goBack)
